# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold plain decimal strings (e.g. "1.002"). Excel would
# otherwise auto-convert such text to a number on assignment, so force Text
# format first to preserve them as the literal strings the source data uses.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '30.191.30'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '1.916.59'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '0.8025'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '244.15'
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '0.3240'
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('D9').Value = '26.79'
$ws.Range('E9').Value = '  +2.14%  '
$ws.Range('D10').Value = '0.07199'
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('D11').Value = '0.7836'
$ws.Range('E11').Value = '  +6.57%  '
$ws.Range('D12').Value = '0.08078'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').Value = '1.922.85'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '5.389'
$ws.Range('E14').Value = '  +4.08%  '
$ws.Range('D15').Value = '93.56'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '30.200.93'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '14.19'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '6.041'
$ws.Range('E18').Value = '  +3.08%  '
$ws.Range('D19').Value = '248.52'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').Value = '0.000007835'
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = '2.173.91'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '8.184'
$ws.Range('E23').Value = '  +18.87%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '0.1630'
$ws.Range('E25').Value = '  +14.94%  '
$ws.Range('D26').Value = '9.451'
$ws.Range('E26').Value = '  +2.82%  '
$ws.Range('D27').Value = '167.20'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').Value = '18.97'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').Value = '2.144'
$ws.Range('E29').Value = '  +5.90%  '
$ws.Range('D30').Value = '1.389'
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').Value = '1.550'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').Value = '4.526'
$ws.Range('E32').Value = '  +5.37%  '
$ws.Range('D33').Value = '0.05737'
$ws.Range('E33').Value = '  +4.13%  '
$ws.Range('D34').Value = '4.144'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('D35').Value = '1.290'
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('D36').Value = '0.7487'
$ws.Range('E36').Value = '  +2.39%  '
$ws.Range('D37').Value = '0.9999'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '2.733'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = '0.01956'
$ws.Range('E39').Value = '  +1.61%  '
$ws.Range('D40').Value = '2.811'
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('D41').Value = '0.4505'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('D42').Value = '73.76'
$ws.Range('E42').Value = '  +2.28%  '
$ws.Range('D43').Value = '6.011'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').Value = '0.8552'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').Value = '1.920'
$ws.Range('E45').Value = '  +2.69%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').Value = '1.036.83'
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('D48').Value = '102.93'
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').Value = '10.05'
$ws.Range('E49').Value = '  +3.40%  '
$ws.Range('D50').Value = '3.098'
$ws.Range('E50').Value = '  +12.22%  '
$ws.Range('D51').Value = '7.613'
$ws.Range('E51').Value = '  +0.86%  '
